$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $text) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue 'D2' '63.643.64'
Set-TextValue 'E2' '  +2.71%  '
Set-TextValue 'D3' '3.129.86'
Set-TextValue 'E3' '  +1.48%  '
Set-TextValue 'D5' '589.05'
Set-TextValue 'E5' '  +1.47%  '
Set-TextValue 'D6' '146.53'
Set-TextValue 'E6' '  +3.08%  '
Set-TextValue 'E7' '  -0.01%  '
Set-TextValue 'D8' '3.122.74'
Set-TextValue 'E8' '  +1.60%  '
Set-TextValue 'E9' '  +1.13%  '
Set-TextValue 'D10' '0.161'
Set-TextValue 'E10' '  +14.60%  '
Set-TextValue 'D11' '5.71'
Set-TextValue 'E11' '  +0.12%  '
Set-TextValue 'E12' '  +0.60%  '
Set-TextValue 'D13' '0.0000252'
Set-TextValue 'E13' '  +4.58%  '
Set-TextValue 'D14' '36.85'
Set-TextValue 'E14' '  +4.27%  '
Set-TextValue 'E15' '  -0.52%  '
Set-TextValue 'D16' '3.642.16'
Set-TextValue 'E16' '  +1.44%  '
Set-TextValue 'D17' '7.17'
Set-TextValue 'E17' '  -1.10%  '
Set-TextValue 'D18' '63.536.70'
Set-TextValue 'E18' '  +2.67%  '
Set-TextValue 'D19' '3.123.94'
Set-TextValue 'E19' '  +1.43%  '
Set-TextValue 'D20' '465.15'
Set-TextValue 'E20' '  +3.96%  '
Set-TextValue 'D21' '14.41'
Set-TextValue 'E21' '  +3.39%  '
Set-TextValue 'E22' '  +0.20%  '
Set-TextValue 'D23' '7.55'
Set-TextValue 'E23' '  +1.36%  '
Set-TextValue 'D24' '13.26'
Set-TextValue 'E24' '  -3.84%  '
Set-TextValue 'D25' '82.20'
Set-TextValue 'E25' '  +0.30%  '
Set-TextValue 'E26' '  -0.11%  '
Set-TextValue 'D27' '8.98'
Set-TextValue 'E27' '  +8.60%  '
Set-TextValue 'D28' '2.71'
Set-TextValue 'E28' '  +1.34%  '
Set-TextValue 'E29' '  -1.79%  '
Set-TextValue 'E30' '  -0.15%  '
Set-TextValue 'E31' '  +1.11%  '
Set-TextValue 'D32' '27.11'
Set-TextValue 'E32' '  +0.92%  '
Set-TextValue 'E33' '  -1.75%  '
Set-TextValue 'D34' '0.0₃0881'
Set-TextValue 'E34' '  +10.94%  '
Set-TextValue 'B35' 'Stacks'
Set-TextValue 'C35' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D35' '2.37'
Set-TextValue 'E35' '  +8.02%  '
Set-TextValue 'B36' 'Mantle'
Set-TextValue 'C36' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D36' '1.05'
Set-TextValue 'E36' '  +1.71%  '
Set-TextValue 'D37' '3.39'
Set-TextValue 'E37' '  +13.48%  '
Set-TextValue 'D38' '6.13'
Set-TextValue 'E38' '  +1.05%  '
Set-TextValue 'D39' '51.05'
Set-TextValue 'E39' '  +1.66%  '
Set-TextValue 'D40' '452.62'
Set-TextValue 'E40' '  +7.24%  '
Set-TextValue 'D41' '8.77'
Set-TextValue 'E41' '  -0.31%  '
Set-TextValue 'D42' '0.0373'
Set-TextValue 'E42' '  +0.49%  '
Set-TextValue 'D43' '2.908.26'
Set-TextValue 'E43' '  +0.32%  '
Set-TextValue 'E44' '  +2.63%  '
Set-TextValue 'E45' '  +1.58%  '
Set-TextValue 'D46' '2.19'
Set-TextValue 'E46' '  +2.39%  '
Set-TextValue 'D47' '36.39'
Set-TextValue 'E47' '  +3.53%  '
Set-TextValue 'D48' '126.14'
Set-TextValue 'E48' '  +1.90%  '
Set-TextValue 'E50' '  +0.01%  '
Set-TextValue 'D51' '24.82'
Set-TextValue 'E51' '  +2.68%  '
